$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.906
$ws.Range("B9").Value = 5.505999999999999
$ws.Range("C9").Value = -11.18
$ws.Range("D9").Value = -7.159000000000001
$ws.Range("B18").Value = 5.275
$ws.Range("B20").Value = 7.069
$ws.Range("C23").Value = -12.953
$ws.Range("C24").Value = -12.548
$ws.Range("C26").Value = -12.401
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("D32").Value = -7.386
$ws.Range("C34").Value = -12.196
$ws.Range("C35").Value = -12.415
$ws.Range("D38").Value = -7.963999999999999
$ws.Range("D45").Value = -7.458999999999999
$ws.Range("C48").Value = -11.529
$ws.Range("D51").Value = -8.059999999999999
$ws.Range("C52").Value = -11.743
$ws.Range("D57").Value = -8.145999999999999
$ws.Range("D64").Value = -7.711
$ws.Range("C66").Value = -11.574
$ws.Range("C67").Value = -10.977
$ws.Range("B69").Value = 5.930999999999999
$ws.Range("B76").Value = 6.248
$ws.Range("C80").Value = -12.372
$ws.Range("B82").Value = 5.496
$ws.Range("D93").Value = -6.981
$ws.Range("C99").Value = -11.793
